# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) / "Valor Mora" (col F) table on rows 16-81
# was re-sorted from descending period order (2112 .. 1607) to ascending
# period order (1607 .. 2112). The mora amounts in column F travel with
# their row, so re-sorting col E in turn reverses the column F sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 16
$endRow = 81

# --- Build the ascending list of periods (YYMM), 1607 .. 2112 ---------
$periods = @()
for ($y = 16; $y -le 21; $y++) {
    for ($m = 1; $m -le 12; $m++) {
        if ($y -eq 16 -and $m -lt 7) { continue }
        $periods += ("{0:D2}{1:D2}" -f $y, $m)
    }
}

# --- Read the current (descending-period-ordered) mora amounts --------
$moraValues = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $moraValues += $ws.Cells.Item($r, 6).Value2
}

# --- Reverse them so they line up with the new ascending period order -
$count = $moraValues.Length
$moraValuesAsc = @()
for ($i = $count - 1; $i -ge 0; $i--) {
    $moraValuesAsc += $moraValues[$i]
}

# --- Write both columns back out ---------------------------------------
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $moraValuesAsc[$i]
}
